# Prediction_template.xlsx edit script
# - Fill in predicted results/goal differences on "Match predictions"
# - Fill in the top goal scorer leaderboard on "Top Goal Scorer"
# - Remove the "about" sheet

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------
# 1. "Match predictions" sheet: predicted result + predicted goal diff
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Match predictions")

# Column E ("Predicted Result") is stored as text (e.g. "2-0"), so force
# a text number format before writing the values.
$ws1.Range("E1:E37").NumberFormat = "@"

$matchResults = @(
    @(2, "2-0", 2),
    @(3, "1-1", 0),
    @(4, "1-0", 1),
    @(5, "2-0", 2),
    @(6, "1-1", 0),
    @(7, "0-1", 1),
    @(8, "0-1", 1),
    @(9, "0-1", 1),
    @(10, "2-0", 2),
    @(11, "0-2", 2),
    @(12, "1-1", 0),
    @(13, "1-0", 1),
    @(14, "2-0", 2),
    @(15, "2-1", 1),
    @(16, "1-1", 0),
    @(17, "1-1", 0),
    @(18, "0-1", 1),
    @(19, "1-0", 1),
    @(20, "0-1", 1),
    @(21, "1-1", 0),
    @(22, "1-1", 0),
    @(23, "0-1", 1),
    @(24, "0-2", 2),
    @(25, "0-2", 2),
    @(26, "1-2", 1),
    @(27, "1-1", 0),
    @(28, "0-3", 3),
    @(29, "0-1", 1),
    @(30, "2-0", 2),
    @(31, "2-0", 2),
    @(32, "2-0", 2),
    @(33, "1-1", 0),
    @(34, "1-1", 0),
    @(35, "0-1", 1),
    @(36, "0-3", 3),
    @(37, "1-1", 0)
)

foreach ($entry in $matchResults) {
    $row = $entry[0]
    $result = $entry[1]
    $goalDiff = $entry[2]
    $ws1.Cells.Item($row, 5).Value = $result
    $ws1.Cells.Item($row, 6).Value = $goalDiff
}

$ws1.Range("F34").Select()

# ---------------------------------------------------------------------
# 2. "Top Goal Scorer" sheet: add Goals column + leaderboard rows
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top Goal Scorer")

$ws2.Cells.Item(1, 4).Value = "Goals"

$topScorers = @(
    @(2, "Romelu Lukaku", "Belgium", 10, 3),
    @(3, "Gonçalo Ramos", "Portugal", 9, 2),
    @(4, "Zeki Amdouni", "Switzerland", 25, 2),
    @(5, "Olivier Giroud", "France", 9, 2),
    @(6, "Cristiano Ronaldo (captain)", "Portugal", 7, 2),
    @(7, "Memphis Depay", "Netherlands", 9, 2),
    @(8, "Harry Kane (captain)", "England", 9, 2),
    @(9, "Niclas Füllkrug", "Germany", 14, 2),
    @(10, "Kylian Mbappé (captain)", "France", 10, 2),
    @(11, "Warren Zaïre-Emery", "France", 33, 2)
)

foreach ($entry in $topScorers) {
    $row = $entry[0]
    $name = $entry[1]
    $country = $entry[2]
    $jersey = $entry[3]
    $goals = $entry[4]
    $ws2.Cells.Item($row, 1).Value = $name
    $ws2.Cells.Item($row, 2).Value = $country
    $ws2.Cells.Item($row, 3).Value = $jersey
    $ws2.Cells.Item($row, 4).Value = $goals
}

$ws2.Range("E13").Select()

# ---------------------------------------------------------------------
# 3. Remove the "about" sheet entirely
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("about")
$ws3.Delete()

# Make "Match predictions" the active sheet/tab when the file is opened.
$ws1.Activate()
$ws1.Range("A5").Select()
